$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.003490139936268095
$ws.Range("C2").Value = 1.531511157193467
$ws.Range("D2").Value = 8.565501471244708
$ws.Range("E2").Value = 2.92668779873165
$ws.Range("F2").Value = 2.995558334026711
$ws.Range("G2").Value = 22

# Row 3
$ws.Range("B3").Value = 0.6398258273924486
$ws.Range("C3").Value = 2.239201287109414
$ws.Range("D3").Value = 20.00114161260279
$ws.Range("E3").Value = 4.472263589347434
$ws.Range("F3").Value = 4.535565317506141
$ws.Range("G3").Value = 21

# Row 4
$ws.Range("B4").Value = -0.9476928289197512
$ws.Range("C4").Value = 1.43433119117039
$ws.Range("D4").Value = 8.057575648683578
$ws.Range("E4").Value = 2.838586910539041
$ws.Range("F4").Value = 2.745226277836465
$ws.Range("G4").Value = 20

# Row 5
$ws.Range("B5").Value = 0.07920638781101008
$ws.Range("C5").Value = 0.4982538277252117
$ws.Range("D5").Value = 0.5017663139879305
$ws.Range("E5").Value = 0.7083546526902541
$ws.Range("F5").Value = 0.7232012390998245
$ws.Range("G5").Value = 19

# Row 6
$ws.Range("B6").Value = 0.1290760995305954
$ws.Range("C6").Value = 0.8563602985157541
$ws.Range("D6").Value = 1.881528676774901
$ws.Range("E6").Value = 1.371688257868712
$ws.Range("F6").Value = 1.405192569417563
$ws.Range("G6").Value = 18

# Row 7
$ws.Range("B7").Value = -0.09004136209350076
$ws.Range("C7").Value = 0.6953463111079407
$ws.Range("D7").Value = 1.386814316799254
$ws.Range("E7").Value = 1.177630806661941
$ws.Range("F7").Value = 1.210320639037886
$ws.Range("G7").Value = 17

# Row 8
$ws.Range("B8").Value = 0.06201971816703138
$ws.Range("C8").Value = 0.466676316002007
$ws.Range("D8").Value = 0.3920755529071259
$ws.Range("E8").Value = 0.6261593670201907
$ws.Range("F8").Value = 0.6449494228223369
$ws.Range("G8").Value = 15
